# Applies the "Added 2021 First Semester Poverty Statistics" commit.
#
# Summary of changes:
#  1. "Full Year" sheet: the Series/year labels in column C shift forward
#     one period (1991->2000, 2000->2006, 2006->2015, 2015->2018) for the
#     existing 52 data rows.
#  2. "First Semester" sheet: the old "Old Series"/"New Series" text labels
#     in column C are replaced with concrete year numbers (2015 / 2018),
#     and 13 new rows of 2021 First-Semester poverty statistics are
#     appended (columns A-C plus H/I).
#  3. The active/selected worksheet moves from "metadata" to
#     "First Semester".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Full Year" sheet -- shift the Series year labels forward one period
# ---------------------------------------------------------------------
$wsFull = $wb.Worksheets.Item("Full Year")

for ($r = 2; $r -le 14; $r++) {
    $wsFull.Range("C$r").Value = 2000
}
for ($r = 15; $r -le 27; $r++) {
    $wsFull.Range("C$r").Value = 2006
}
for ($r = 28; $r -le 40; $r++) {
    $wsFull.Range("C$r").Value = 2015
}
for ($r = 41; $r -le 53; $r++) {
    $wsFull.Range("C$r").Value = 2018
}

# ---------------------------------------------------------------------
# 2. "First Semester" sheet
# ---------------------------------------------------------------------
$wsFS = $wb.Worksheets.Item("First Semester")

# Replace "Old Series" text labels (rows 2-14) with the literal year 2015
for ($r = 2; $r -le 14; $r++) {
    $wsFS.Range("C$r").Value = 2015
}

# Replace "New Series" text labels (rows 15-27) with the literal year 2018
for ($r = 15; $r -le 27; $r++) {
    $wsFS.Range("C$r").Value = 2018
}

# Append the new 2021 First Semester data block (rows 28-40)
$newRows = @(
    @{ Row=28; A="Annual Per Capita Poverty Threshold";       B="pesos";             H=12637.918553006057;  I=14498.078114544631 },
    @{ Row=29; A="Poverty Incidence among Families";          B="percent";           H=16.187536428016173;  I=17.955553145557833 },
    @{ Row=30; A="Magnitude of Poor Families";                B="households ('000)"; H=4039.4099904000514;  I=4739.8110360999972 },
    @{ Row=31; A="Poverty Incidence among Population";        B="percent";           H=21.050909176398534;  I=23.717926630826302 },
    @{ Row=32; A="Magnitude of Poor Population";               B="persons ('000)";    H=22262.427214906216;  I=26136.836813200051 },
    @{ Row=33; A="Annual Per Capita Food Threshold";          B="pesos";             H=8848.5788319105577;  I=10071.257720130052 },
    @{ Row=34; A="Subsistence Incidence among Families";      B="percent";           H=6.2078326218868618;  I=7.0792336675070695 },
    @{ Row=35; A="Magnitude of Subsistence Poor Families";    B="households ('000)"; H=1549.0918721999913;  I=1868.7383002000088 },
    @{ Row=36; A="Subsistence Incidence among Population";    B="percent";           H=8.5392721451641069;  I=9.9320081650149685 },
    @{ Row=37; A="Magnitude of Subsistence Poor Population";  B="persons ('000)";    H=9030.7227591447827;  I=10944.939693799992 },
    @{ Row=38; A="Income Gap";                                 B="percent";           H=26.97636346869961;   I=27.037492166175564 },
    @{ Row=39; A="Poverty Gap";                                B="percent";           H=4.3668086634497962;  I=4.8547312751236902 },
    @{ Row=40; A="Severity of Poverty";                        B="percent";           H=1.7974468063177669;  I=1.9483855219994111 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $wsFS.Range("A$r").Value = $nr.A
    $wsFS.Range("B$r").Value = $nr.B
    $wsFS.Range("C$r").Value = 2021
    $wsFS.Range("H$r").Value = $nr.H
    $wsFS.Range("H$r").NumberFormat = "#,##0.00"
    $wsFS.Range("I$r").Value = $nr.I
    $wsFS.Range("I$r").NumberFormat = "#,##0.00"
}

# ---------------------------------------------------------------------
# 3. Make "First Semester" the active / selected sheet
# ---------------------------------------------------------------------
$wsFS.Activate()
$wsFS.Select()
